$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.939.81'
$ws.Range("E2").Value = '  +2.73%  '
$ws.Range("D3").Value = '3.062.40'
$ws.Range("E3").Value = '  +3.03%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '526.22'
$ws.Range("E5").Value = '  +6.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.39'
$ws.Range("E6").Value = '  +5.63%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.447'
$ws.Range("E8").Value = '  +5.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.67'
$ws.Range("E9").Value = '  +6.92%  '
$ws.Range("E10").Value = '  +7.75%  '
$ws.Range("E11").Value = '  +5.67%  '
$ws.Range("E12").Value = '  +2.49%  '
$ws.Range("D13").Value = '3.589.83'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.23'
$ws.Range("E14").Value = '  +8.47%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000170'
$ws.Range("E15").Value = '  +16.52%  '
$ws.Range("D16").Value = '57.903.25'
$ws.Range("E16").Value = '  +2.72%  '
$ws.Range("E17").Value = '  +7.64%  '
$ws.Range("D18").Value = '3.063.05'
$ws.Range("E18").Value = '  +3.25%  '
$ws.Range("E19").Value = '  +6.39%  '
$ws.Range("E20").Value = '  +5.21%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '339.89'
$ws.Range("E21").Value = '  +4.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("E23").Value = '  -0.67%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.503'
$ws.Range("E24").Value = '  +7.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.00'
$ws.Range("E25").Value = '  +5.79%  '
$ws.Range("E26").Value = '  +6.80%  '
$ws.Range("D27").Value = '0.0₃0976'
$ws.Range("E27").Value = '  +8.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.97'
$ws.Range("E29").Value = '  +7.12%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.36'
$ws.Range("E30").Value = '  +9.11%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.85'
$ws.Range("E31").Value = '  +6.60%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.24'
$ws.Range("E32").Value = '  +6.29%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.15'
$ws.Range("E33").Value = '  +4.77%  '
$ws.Range("E34").Value = '  +7.38%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '156.65'
$ws.Range("E35").Value = '  +2.50%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.01'
$ws.Range("E36").Value = '  +7.22%  '
$ws.Range("E37").Value = '  +4.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '26.37'
$ws.Range("E38").Value = '  +14.21%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0703'
$ws.Range("E39").Value = '  +4.56%  '
$ws.Range("D40").Value = '3.096.20'
$ws.Range("E40").Value = '  +3.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '37.89'
$ws.Range("E41").Value = '  +3.52%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.91'
$ws.Range("E42").Value = '  +9.84%  '
$ws.Range("E43").Value = '  +0.26%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.665'
$ws.Range("E44").Value = '  +4.40%  '
$ws.Range("E45").Value = '  +5.96%  '
$ws.Range("D46").Value = '2.333.34'
$ws.Range("E46").Value = '  +5.69%  '
$ws.Range("E47").Value = '  +3.58%  '
$ws.Range("E48").Value = '  +3.55%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0245'
$ws.Range("E49").Value = '  +3.53%  '
$ws.Range("E50").Value = '  +5.30%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.26'
$ws.Range("E51").Value = '  +6.46%  '
